$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Terms")

# Column A (rows 2-9) previously held "Prog" as a short module label;
# replace it with the full module title used elsewhere in the workbook.
$ws.Range("A2:A9").Value = "Introduction to Programming"

# Give column A an explicit best-fit width now that it holds longer text.
$ws.Columns.Item(1).ColumnWidth = 20.28

# Update the saved selection to match the new content in column A.
$ws.Range("A2:A9").Select()
